$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap column pair (B,C) with column pair (D,E) for rows 1-19, then fix
# number formats: the "Start Time" values (originally in B, with an
# h:mm:ss style) move to D, and B must lose that time formatting.
for ($r = 1; $r -le 19; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value()
    $cVal = $ws.Cells.Item($r, 3).Value()
    $dVal = $ws.Cells.Item($r, 4).Value()
    $eVal = $ws.Cells.Item($r, 5).Value()

    $ws.Cells.Item($r, 2).Value = $dVal
    $ws.Cells.Item($r, 3).Value = $eVal
    $ws.Cells.Item($r, 4).Value = $bVal
    $ws.Cells.Item($r, 5).Value = $cVal
}

# Column B no longer holds time values; reset its formatting back to the
# workbook default (removes the inherited h:mm:ss style cleanly).
$ws.Range("B1:B19").Style = "Normal"

# Column D now holds the time-of-day values; give it the same h:mm:ss
# number format that column B used to carry.
$ws.Range("D2:D19").NumberFormat = "h:mm:ss"

# Update the remembered selection as recorded in the saved workbook.
$ws.Range("D24").Select()
